$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'228.58"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.Value = "'22.33"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.Value = "'5.272"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.Value = "'0.05536"
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.Value = "'3.387"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.Value = "'6.467"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.Value = "'1.076"
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.Value = "'0.7709"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'0.07422"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.Value = "'0.03157"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.Value = "'0.02947"
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.Value = "'0.09265"
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'0.001660"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.Value = "'3.261"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.Value = "'0.04775"
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.Value = "'0.01166"
$cell.Style = "Normal"

$ws.Range("E18").Value = "17OneONEBestin24h"

$cell = $ws.Range("D19")
$cell.Value = "'0.006206"
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.Value = "'0.005226"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.Value = "'0.001064"
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.Value = "'0.0001502"
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.Value = "'3.907"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.Value = "'0.1285"
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.Value = "'0.0005005"
$cell.Style = "Normal"

$ws.Range("E27").Value = "26UpBotsUBXT"

$cell = $ws.Range("D40")
$cell.Value = "'0.03944"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.Value = "'0.007123"
$cell.Style = "Normal"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$cell = $ws.Range("D42")
$cell.Value = "'0.003504"
$cell.Style = "Normal"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$cell = $ws.Range("D43")
$cell.Value = "'0.1033"
$cell.Style = "Normal"

$ws.Range("E43").Value = "42BKEXTokenBKK"

$cell = $ws.Range("D44")
$cell.Value = "'0.008754"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.Value = "'0.00005442"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.Value = "'0.7860"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.Value = "'0.04102"
$cell.Style = "Normal"

$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
